$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.292711734771729
$ws.Range("B1").Value = 2.941659927368164
$ws.Range("C1").Value = 5.214778423309326
$ws.Range("D1").Value = 1.842502117156982
$ws.Range("E1").Value = 1.012032866477966
